# Auto-generated edit script: applies scheduled market-data refresh updates
# to the Leve profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H123").Value = 42110
$ws.Range("J123").Value = 42110
$ws.Range("L123").Value = 42110
$ws.Range("N123").Value = -51910

$ws.Range("H131").Value = 2069.25
$ws.Range("J131").Value = 3500.2
$ws.Range("L131").Value = 10500.6
$ws.Range("N131").Value = -20580.6

$ws.Range("H132").Value = 1071.3823
$ws.Range("I132").Value = 1091.2258
$ws.Range("K132").Value = 3273.6774
$ws.Range("M132").Value = -743.6773999999996

$ws.Range("H135").Value = 716.5
$ws.Range("I135").Value = 716.5
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 6448.5
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -3913.5
$ws.Range("N135").ClearContents()

$ws.Range("H141").Value = 3535.4443
$ws.Range("I141").Value = 1460
$ws.Range("K141").Value = 4380
$ws.Range("M141").Value = 800

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2772.423
$ws.Range("I61").Value = 2149.7273
$ws.Range("J61").Value = 6197.25
$ws.Range("K61").Value = 2149.7273
$ws.Range("L61").Value = 6197.25
$ws.Range("M61").Value = -1937.7273
$ws.Range("N61").Value = -6621.25

$ws.Range("H136").Value = 2772.423
$ws.Range("I136").Value = 2149.7273
$ws.Range("J136").Value = 6197.25
$ws.Range("K136").Value = 6449.1819
$ws.Range("L136").Value = 18591.75
$ws.Range("M136").Value = -3899.1819
$ws.Range("N136").Value = -23691.75

$ws.Range("H139").Value = 49143.332
$ws.Range("J139").Value = 49143.332
$ws.Range("L139").Value = 49143.332
$ws.Range("N139").Value = -59423.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1795.1177
$ws.Range("I107").Value = 1955.4445
$ws.Range("J107").Value = 1614.75
$ws.Range("K107").Value = 1955.4445
$ws.Range("L107").Value = 1614.75
$ws.Range("M107").Value = -35.44450000000006
$ws.Range("N107").Value = -5454.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 59806
$ws.Range("I23").Value = 20000
$ws.Range("J23").Value = 69757.5
$ws.Range("K23").Value = 20000
$ws.Range("L23").Value = 69757.5
$ws.Range("M23").Value = -19760
$ws.Range("N23").Value = -70237.5

$ws.Range("H27").Value = 59806
$ws.Range("I27").Value = 20000
$ws.Range("J27").Value = 69757.5
$ws.Range("K27").Value = 20000
$ws.Range("L27").Value = 69757.5
$ws.Range("M27").Value = -19808
$ws.Range("N27").Value = -70141.5

$ws.Range("H31").Value = 4089.3076
$ws.Range("I31").Value = 1582.625
$ws.Range("J31").Value = 8100
$ws.Range("K31").Value = 1582.625
$ws.Range("L31").Value = 8100
$ws.Range("M31").Value = -1287.625
$ws.Range("N31").Value = -8690

$ws.Range("H34").Value = 4089.3076
$ws.Range("I34").Value = 1582.625
$ws.Range("J34").Value = 8100
$ws.Range("K34").Value = 1582.625
$ws.Range("L34").Value = 8100
$ws.Range("M34").Value = -1380.625
$ws.Range("N34").Value = -8504

$ws.Range("H132").Value = 2126.8948
$ws.Range("I132").Value = 1517.5
$ws.Range("K132").Value = 4552.5
$ws.Range("M132").Value = -2022.5

$ws.Range("H134").Value = 847.2941
$ws.Range("I134").Value = 846.8461
$ws.Range("J134").Value = 848.75
$ws.Range("K134").Value = 2540.5383
$ws.Range("L134").Value = 2546.25
$ws.Range("M134").Value = -5.538300000000163
$ws.Range("N134").Value = -7616.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 769.5625
$ws.Range("J107").Value = 800.6667
$ws.Range("L107").Value = 2402.0001
$ws.Range("N107").Value = -6242.0001

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 10000
$ws.Range("J44").Value = 10000
$ws.Range("L44").Value = 10000
$ws.Range("N44").Value = -11192

$ws.Range("H80").Value = 2440
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 2440
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws.Range("H61").Value = 2174.9473
$ws.Range("I61").Value = 2191.4
$ws.Range("J61").Value = 2113.25
$ws.Range("K61").Value = 2191.4
$ws.Range("L61").Value = 2113.25
$ws.Range("M61").Value = -1989.4
$ws.Range("N61").Value = -2517.25

$ws.Range("H113").Value = 2174.9473
$ws.Range("I113").Value = 2191.4
$ws.Range("J113").Value = 2113.25
$ws.Range("K113").Value = 2191.4
$ws.Range("L113").Value = 2113.25
$ws.Range("M113").Value = -21.40000000000009
$ws.Range("N113").Value = -6453.25

$ws.Range("H136").Value = 2750.2222
$ws.Range("I136").Value = 1471.6364
$ws.Range("J136").Value = 4759.4287
$ws.Range("K136").Value = 4414.9092
$ws.Range("L136").Value = 14278.2861
$ws.Range("M136").Value = -1864.9092
$ws.Range("N136").Value = -19378.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 2266.6667
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 2266.6667
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 2266.6667
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -2724.6667

$ws.Range("H126").Value = 4242.7026
$ws.Range("I126").Value = 4076.8708
$ws.Range("J126").Value = 5099.5
$ws.Range("K126").Value = 12230.6124
$ws.Range("L126").Value = 15298.5
$ws.Range("M126").Value = -9760.6124
$ws.Range("N126").Value = -20238.5

$ws.Range("H132").Value = 5233.628
$ws.Range("I132").Value = 1136.7778
$ws.Range("K132").Value = 3410.3334
$ws.Range("M132").Value = -880.3334000000004

$ws.Range("H136").Value = 16342000
$ws.Range("I136").Value = 25254466
$ws.Range("J136").Value = 2481.3333
$ws.Range("K136").Value = 75763398
$ws.Range("L136").Value = 7443.999899999999
$ws.Range("M136").Value = -75760848
$ws.Range("N136").Value = -12543.9999
